# Update column G ("K" = strikeouts) values on Sheet1, rows 2-31,
# as part of regenerating save_data to use K instead of Strike#.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 3
    3  = 8
    4  = 7
    5  = 5
    6  = 3
    7  = 4
    8  = 6
    9  = 9
    10 = 7
    11 = 0
    12 = 3
    13 = 4
    14 = 4
    15 = 6
    16 = 5
    17 = 8
    18 = 7
    19 = 4
    20 = 7
    21 = 6
    22 = 5
    23 = 4
    24 = 1
    25 = 5
    26 = 3
    27 = 0
    28 = 6
    29 = 2
    30 = 3
    31 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
